# Adigeni "Number of stillbirths" sheet - upgrade left table (commit: "upgrade
# left table until javakheti"): rename the sheet, mark the Urban row as fully
# confidential, normalise a few more Rural cells to confidential, replace the
# single-glyph ellipsis ("…") with a plain three-dot ellipsis ("...") used
# throughout the table, and drop the spare blank row above the footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from the generic "1" to "Adigeni".
$ws.Name = "Adigeni"

$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")

# Row 6 ("Urban"): every year becomes confidential/unavailable.
foreach ($col in $dataCols) {
    $ws.Range($col + "6").Value = "..."
}

# Row 7 ("Rural"): 2013, 2017 and 2020 newly become confidential/unavailable;
# the remaining years keep their existing values.
foreach ($col in @("E","I","L")) {
    $ws.Range($col + "7").Value = "..."
}

# The table already used an ellipsis marker for confidential/unavailable
# cells ("…"); normalise every such marker still left in rows 5 and 7 to the
# plain three-dot form ("...") used elsewhere in the sheet.
foreach ($col in @("C","D","F","H","K","M","N","O")) {
    $ws.Range($col + "5").Value = "..."
    $ws.Range($col + "7").Value = "..."
}

# Remove the spare blank row 8 so the footnote (previously row 9) moves up to
# row 8, matching the tightened table layout.
$ws.Rows("8").Delete()
